# Update "want to go" counts (column F) on the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1834
$ws1.Range("F5").Value = 796
$ws1.Range("F10").Value = 238
$ws1.Range("F13").Value = 133
$ws1.Range("F15").Value = 4299
$ws1.Range("F17").Value = 31
$ws1.Range("F20").Value = 986
$ws1.Range("F21").Value = 1569
$ws1.Range("F22").Value = 362
$ws1.Range("F26").Value = 2011
$ws1.Range("F31").Value = 57

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1834
$ws4.Range("F5").Value = 796
$ws4.Range("F10").Value = 238
$ws4.Range("F13").Value = 133
$ws4.Range("F15").Value = 4299
$ws4.Range("F17").Value = 31
$ws4.Range("F20").Value = 986
$ws4.Range("F21").Value = 1570
$ws4.Range("F22").Value = 362
$ws4.Range("F26").Value = 2011
$ws4.Range("F31").Value = 57
